$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds the last-changed date for each record.
# Rows 2-56 all move from serial date 45233 (2023-11-03) to 45243 (2023-11-13).
$ws.Range("C2:C56").Value = 45243
